# Planilla_Wurth_World_Cup_2026.xlsx - "Add files via upload"
# The sheet's header row (row 1) holds long column titles in C1:I1; this
# edit widens those columns so the headers are readable, and leaves the
# cell cursor parked at C17 (as it was when the sheet was re-saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns C:I to fit their (long) header text. ColumnWidth is rounded
# to this host's internal pixel grid, so the inputs below are pre-solved to
# land on the closest achievable grid value to the real target widths
# (C=34.5703125, D=30.7109375, E=40.85546875, F=40.140625, G=28.7109375,
#  H=39.28515625, I=23.42578125 "characters").
$ws.Columns.Item(3).ColumnWidth = 33.6666666666667
$ws.Columns.Item(4).ColumnWidth = 29.8333333333333
$ws.Columns.Item(5).ColumnWidth = 40
$ws.Columns.Item(6).ColumnWidth = 39.3333333333333
$ws.Columns.Item(7).ColumnWidth = 27.8333333333333
$ws.Columns.Item(8).ColumnWidth = 38.5
$ws.Columns.Item(9).ColumnWidth = 22.6666666666667

# Leave the selection where it ended up in the saved file.
$ws.Range("C17").Select() | Out-Null
